# Replace iCE5LP1K with iCE40UP5K and refresh KiCost pricing/availability data
# for row 20 (U6 - the FPGA), plus refresh the two KiCost run timestamps.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Part number changes (U6: FPGA) ---
$ws.Range("B20").Value = "ICE40UP5K-SG48"
$ws.Range("G20").Value = "ICE40UP5K-SG48I"
$ws.Range("V20").Value = "842-ICE40UP5K-SG48I"

# --- Digi-Key (K..P20) pricing/availability refresh ---
# Digi-Key catalog # changed along with the part switch.
$ws.Range("P20").Value = "220-2212-1-ND"

# Digi-Key no longer stocks this exact catalog number - mark K20 as "NonStk"
# (copy the grayed-out "not stocked" number format from another NonStk cell,
# then overwrite the value with the NonStk text).
$ws.Range("Q9").Copy()
$ws.Range("K20").PasteSpecial(-4122)
$ws.Range("K20").Value = "NonStk"

# Updated Digi-Key price breaks: 1/25/100 qty -> $10.10/$8.80/$8.45
$ws.Range("M20").Formula = '=IFERROR(IF(OR(L20>=N20,H20>=N20),LOOKUP(IF(L20="",H20,L20),{0,1,25,100},{0.0,10.1,8.8,8.45}),"MOQ="&N20),"")'

# --- Mouser (Q..V20) pricing/availability refresh ---
$ws.Range("Q20").Value = 12647

# Updated Mouser price breaks: 1/25/100/2000 qty -> $10.17/$8.85/$8.50/$8.50
$ws.Range("S20").Formula = '=IFERROR(IF(OR(R20>=T20,H20>=T20),LOOKUP(IF(R20="",H20,R20),{0,1,25,100,2000},{0.0,10.17,8.85,8.5,8.5}),"MOQ="&T20),"")'

# --- Comments ---
$null = $ws.Range("K20").AddComment("This part is listed but is not stocked.")

$digikeyBreaks = @'
Qty/Price Breaks (USD):
  Qty  -  Unit$  -  Ext$
================
     1  $10.10     $10.10
    25   $8.80    $220.00
   100   $8.45    $845.00
'@
$null = $ws.Range("M20").Comment.Text($digikeyBreaks)

$null = $ws.Range("Q20").Comment.Text("12647 In Stock")

$mouserBreaks = @'
Qty/Price Breaks (USD):
  Qty  -  Unit$  -  Ext$
================
     1  $10.17     $10.17
    25   $8.85    $221.25
   100   $8.50    $850.00
  2000   $8.50 $17,000.00
'@
$null = $ws.Range("S20").Comment.Text($mouserBreaks)

$null = $ws.Range("V20").Comment.Text("Desc: FPGA - Field Programmable Gate Array iCE40 UltraPlus, 5280 LUTs, 1.2V")

# --- Refresh KiCost run timestamps ---
$ws.Range("B3").Value = "Sat 04 Nov 2023 08:29:48 AM CET"
$ws.Range("B4").Value = "2023-11-04 08:29:49"
